# Insert three new data rows (253-255) into the "Fruta / Palta" weekly
# price table for Terminal Hortofrutícola Agro Chillán, pushing the
# existing rows 253-269 down to 256-272.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows before the current row 253, shifting everything
# below (including formatting) down by 3 rows.
$ws.Rows("253:255").Insert()

# Common (constant-across-block) column values shared by every row in
# this Terminal Hortofrutícola Agro Chillán / Palta subset.
$marketId   = 7
$market     = "Terminal Hortofrutícola Agro Chillán"
$region     = "Ñuble"
$codreg     = 16
$tipo       = "Fruta"
$prodId     = 100106
$producto   = "Oleaginosos"
$catId      = 100106002
$categoria  = "Palta"
$variedad   = "Hass"

function Set-PalataRow($r, $fecha, $calidad, $volumen, $pmin, $pmax, $pprom, $unidad, $origen, $precioKg, $kgUnidad) {
    $ws.Cells.Item($r, 1).Value  = $marketId
    $ws.Cells.Item($r, 2).Value  = $market
    $ws.Cells.Item($r, 3).Value  = $region
    $ws.Cells.Item($r, 4).Value  = $fecha
    $ws.Cells.Item($r, 5).Value  = $codreg
    $ws.Cells.Item($r, 6).Value  = $tipo
    $ws.Cells.Item($r, 7).Value  = $prodId
    $ws.Cells.Item($r, 8).Value  = $producto
    $ws.Cells.Item($r, 9).Value  = $catId
    $ws.Cells.Item($r, 10).Value = $categoria
    $ws.Cells.Item($r, 11).Value = $variedad
    $ws.Cells.Item($r, 12).Value = $calidad
    $ws.Cells.Item($r, 13).Value = $volumen
    $ws.Cells.Item($r, 14).Value = $pmin
    $ws.Cells.Item($r, 15).Value = $pmax
    $ws.Cells.Item($r, 16).Value = $pprom
    $ws.Cells.Item($r, 17).Value = $unidad
    $ws.Cells.Item($r, 18).Value = $origen
    $ws.Cells.Item($r, 19).Value = $precioKg
    $ws.Cells.Item($r, 20).Value = $kgUnidad
}

# Row 253
Set-PalataRow 253 44461 "1a nueva(o)" 120 2800 2900 2850 "$/kilo (en caja de 15 kilos)" "Provincia de Quillota" 2850 1

# Row 254
Set-PalataRow 254 44461 "Primera" 360 25000 26000 25500 "$/bandeja 10 kilos" "Perú" 2550 10

# Row 255
Set-PalataRow 255 44461 "Segunda" 100 24000 24000 24000 "$/bandeja 10 kilos" "Perú" 2400 10
